# "Fixed typos in code in comparison slide"
# Slide 23 ("Assembly Code Examples") has three code-comparison textboxes
# (C/C++, MIPS, ARM). This fixes:
#   - C/C++:  "x < 10=" -> "x <= 10"            (typo'd operator)
#   - MIPS:   "bnez $4, skip" -> "bnez $4, sk"   (label is "sk:", not "skip:")
#   - MIPS:   "j loop"       -> "j lp"           (label is "lp:",  not "loop:")
#   - ARM:    "BGT exit"     -> "BGT ex"         (label is "ex:",  not "exit:")
#   - ARM:    "BAL  loop"    -> "BAL  lp"        (label is "lp:",  not "loop:")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)

# --- C/C++ textbox ---------------------------------------------------------
$ccShape = $s.Shapes.Item(4)
$ccRange = $ccShape.TextFrame.TextRange
$ccOld = "for(int x = 1; x < 10=; x++)"
$ccNew = "for(int x = 1; x <= 10; x++)"
$ccIdx = $ccRange.Text.IndexOf($ccOld)
if ($ccIdx -ge 0) {
    $ccRange.Characters($ccIdx + 1, $ccOld.Length).Text = $ccNew
}

# --- MIPS textbox -----------------------------------------------------------
$mipsShape = $s.Shapes.Item(5)
$mipsRange = $mipsShape.TextFrame.TextRange

# "bnez $4, skip" -> "bnez $4, sk"  (splits the trailing run in two, like
# retyping just the "skip" word down to "sk")
$skipIdx = $mipsRange.Text.IndexOf("skip")
if ($skipIdx -ge 0) {
    $mipsRange.Characters($skipIdx + 1, 4).Text = "sk"
}

# "     j loop" -> "     j lp"
$jIdx = $mipsRange.Text.IndexOf("j loop")
$loopIdx = $mipsRange.Text.IndexOf("loop", $jIdx)
if ($loopIdx -ge 0) {
    $mipsRange.Characters($loopIdx + 1, 4).Text = "lp"
}

# --- ARM textbox -------------------------------------------------------------
$armShape = $s.Shapes.Item(6)
$armRange = $armShape.TextFrame.TextRange

# "    BGT exit" -> "    BGT ex" -- plain in-place text edit, no run split
$bgtOld = "    BGT exit"
$bgtNew = "    BGT ex"
$bgtIdx = $armRange.Text.IndexOf($bgtOld)
if ($bgtIdx -ge 0) {
    $armRange.Characters($bgtIdx + 1, $bgtOld.Length).Text = $bgtNew
}

# "    BAL  loop" -> "    BAL  lp"
$balIdx = $armRange.Text.IndexOf("BAL  loop")
$balLoopIdx = $armRange.Text.IndexOf("loop", $balIdx)
if ($balLoopIdx -ge 0) {
    $armRange.Characters($balLoopIdx + 1, 4).Text = "lp"
}
